$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G gets a literal (non-formula) copy of column F's computed string,
# for every data row (1-50) - mirrors "UPDATE oc_country ..." values already
# produced by the F-column formulas.
for ($r = 1; $r -le 50; $r++) {
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 6).Value2
}

# New column G needs the same custom width treatment as the other text columns.
$ws.Columns.Item(7).ColumnWidth = 68.57142857142857

# Selection/view moves from column F to the newly added column G, and the
# sheet no longer needs to keep B1 pinned as the top-left cell.
[void]$ws.Range("G1:G1048576").Select()
